$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgi2"
$ws.Range("C2").Value = "Adam22"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1188046666666667
$ws.Range("H2").Value = 0.356414
$ws.Range("I2").Value = 0.00979853232878679
$ws.Range("J2").Value = 0.009798532328786792
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.418229333333334
$ws.Range("N2").Value = 19.254688
$ws.Range("O2").Value = 0.4263166911921744
$ws.Range("P2").Value = 0.4263166911921744
$ws.Range("Q2").Value = 0.7625155965368889
$ws.Range("R2").Value = 6.862640368832
$ws.Range("S2").Value = 0.004177277880947935
$ws.Range("T2").Value = 0.004177277880947936

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lgi2"
$ws.Range("C3").Value = "Adam22"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1188046666666667
$ws.Range("H3").Value = 0.356414
$ws.Range("I3").Value = 0.00979853232878679
$ws.Range("J3").Value = 0.009798532328786792
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.860692
$ws.Range("N3").Value = 11.582076
$ws.Range("O3").Value = 0.256437929165941
$ws.Range("P3").Value = 0.256437929165941
$ws.Range("Q3").Value = 0.4586682261626667
$ws.Range("R3").Value = 4.128014035464
$ws.Range("S3").Value = 0.00251271533925961
$ws.Range("T3").Value = 0.00251271533925961

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lgi2"
$ws.Range("C4").Value = "Adam22"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1188046666666667
$ws.Range("H4").Value = 0.356414
$ws.Range("I4").Value = 0.00979853232878679
$ws.Range("J4").Value = 0.009798532328786792
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.776152666666666
$ws.Range("N4").Value = 14.328458
$ws.Range("O4").Value = 0.3172453796418847
$ws.Range("P4").Value = 0.3172453796418846
$ws.Range("Q4").Value = 0.5674292255124445
$ws.Range("R4").Value = 5.106863029612
$ws.Range("S4").Value = 0.003108539108579245
$ws.Range("T4").Value = 0.003108539108579245

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lgi2"
$ws.Range("C5").Value = "Adam22"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.502271
$ws.Range("H5").Value = 34.50681299999999
$ws.Range("I5").Value = 0.9486611714015168
$ws.Range("J5").Value = 0.948661171401517
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.418229333333334
$ws.Range("N5").Value = 19.254688
$ws.Range("O5").Value = 0.4263166911921744
$ws.Range("P5").Value = 0.4263166911921744
$ws.Range("Q5").Value = 73.82421313214932
$ws.Range("R5").Value = 664.4179181893439
$ws.Range("S5").Value = 0.4044300916543869
$ws.Range("T5").Value = 0.4044300916543869

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lgi2"
$ws.Range("C6").Value = "Adam22"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.502271
$ws.Range("H6").Value = 34.50681299999999
$ws.Range("I6").Value = 0.9486611714015168
$ws.Range("J6").Value = 0.948661171401517
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.860692
$ws.Range("N6").Value = 11.582076
$ws.Range("O6").Value = 0.256437929165941
$ws.Range("P6").Value = 0.256437929165941
$ws.Range("Q6").Value = 44.406725631532
$ws.Range("R6").Value = 399.660530683788
$ws.Range("S6").Value = 0.2432727062743408
$ws.Range("T6").Value = 0.2432727062743408

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lgi2"
$ws.Range("C7").Value = "Adam22"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.502271
$ws.Range("H7").Value = 34.50681299999999
$ws.Range("I7").Value = 0.9486611714015168
$ws.Range("J7").Value = 0.948661171401517
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.776152666666666
$ws.Range("N7").Value = 14.328458
$ws.Range("O7").Value = 0.3172453796418847
$ws.Range("P7").Value = 0.3172453796418846
$ws.Range("Q7").Value = 54.93660230937266
$ws.Range("R7").Value = 494.4294207843539
$ws.Range("S7").Value = 0.3009583734727892
$ws.Range("T7").Value = 0.3009583734727892

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Lgi2"
$ws.Range("C8").Value = "Adam22"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5036653333333333
$ws.Range("H8").Value = 1.510996
$ws.Range("I8").Value = 0.04154029626969626
$ws.Range("J8").Value = 0.04154029626969627
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.418229333333334
$ws.Range("N8").Value = 19.254688
$ws.Range("O8").Value = 0.4263166911921744
$ws.Range("P8").Value = 0.4263166911921744
$ws.Range("Q8").Value = 3.232639616583111
$ws.Range("R8").Value = 29.093756549248
$ws.Range("S8").Value = 0.01770932165683953
$ws.Range("T8").Value = 0.01770932165683954

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Lgi2"
$ws.Range("C9").Value = "Adam22"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5036653333333333
$ws.Range("H9").Value = 1.510996
$ws.Range("I9").Value = 0.04154029626969626
$ws.Range("J9").Value = 0.04154029626969627
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.860692
$ws.Range("N9").Value = 11.582076
$ws.Range("O9").Value = 0.256437929165941
$ws.Range("P9").Value = 0.256437929165941
$ws.Range("Q9").Value = 1.944496723077333
$ws.Range("R9").Value = 17.500470507696
$ws.Range("S9").Value = 0.01065250755234057
$ws.Range("T9").Value = 0.01065250755234057

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Lgi2"
$ws.Range("C10").Value = "Adam22"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5036653333333333
$ws.Range("H10").Value = 1.510996
$ws.Range("I10").Value = 0.04154029626969626
$ws.Range("J10").Value = 0.04154029626969627
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.776152666666666
$ws.Range("N10").Value = 14.328458
$ws.Range("O10").Value = 0.3172453796418847
$ws.Range("P10").Value = 0.3172453796418846
$ws.Range("Q10").Value = 2.405582524907555
$ws.Range("R10").Value = 21.650242724168
$ws.Range("S10").Value = 0.01317846706051615
$ws.Range("T10").Value = 0.01317846706051616

